$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.499.75'
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("D3").Value = '1.570.62'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -1.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.492'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -1.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.76'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.250'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0600'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("D12").Value = '1.794.50'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '1.578.33'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '27.485.06'
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '226.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '0.0₃0708'
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.992'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = '1.457.20'
$ws.Range("E33").Value = '  +2.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("E35").Value = '  +3.12%  '
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.991'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.973'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("D47").Value = '1.706.66'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0945'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.93%  '
